$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1160.6792
$ws.Range("I15").Value = 1160.6792
$ws.Range("K15").Value = 3482.0376
$ws.Range("M15").Value = -3313.0376

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 11375.667
$ws.Range("I51").Value = 18666.834
$ws.Range("J51").Value = 4084.5
$ws.Range("K51").Value = 18666.834
$ws.Range("L51").Value = 4084.5
$ws.Range("M51").Value = -18182.834
$ws.Range("N51").Value = -5052.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 11912025
$ws.Range("I132").Value = 13165500
$ws.Range("K132").Value = 39496500
$ws.Range("M132").Value = -39493970

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1827.3112
$ws.Range("I137").Value = 1432.5161
$ws.Range("J137").Value = 2701.5
$ws.Range("K137").Value = 4297.5483
$ws.Range("L137").Value = 8104.5
$ws.Range("M137").Value = -1747.5483
$ws.Range("N137").Value = -13204.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 34998.395
$ws.Range("I32").Value = 7000.1816
$ws.Range("J32").Value = 81195.45
$ws.Range("K32").Value = 7000.1816
$ws.Range("L32").Value = 81195.45
$ws.Range("M32").Value = -6713.1816
$ws.Range("N32").Value = -81769.45

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2317.875
$ws.Range("I61").Value = 1295
$ws.Range("J61").Value = 3048.5
$ws.Range("K61").Value = 1295
$ws.Range("L61").Value = 3048.5
$ws.Range("M61").Value = -1083
$ws.Range("N61").Value = -3472.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1901.75
$ws.Range("I74").Value = 1505.25
$ws.Range("J74").Value = 2694.75
$ws.Range("K74").Value = 1505.25
$ws.Range("L74").Value = 2694.75
$ws.Range("M74").Value = -631.25
$ws.Range("N74").Value = -4442.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1901.75
$ws.Range("I77").Value = 1505.25
$ws.Range("J77").Value = 2694.75
$ws.Range("K77").Value = 7526.25
$ws.Range("L77").Value = 13473.75
$ws.Range("M77").Value = -3158.25
$ws.Range("N77").Value = -22209.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H119").Value = 40644
$ws.Range("J119").Value = 40644
$ws.Range("L119").Value = 40644
$ws.Range("N119").Value = -50320

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2745.1
$ws.Range("I132").Value = 2974.3333
$ws.Range("J132").Value = 2210.2222
$ws.Range("K132").Value = 8922.999899999999
$ws.Range("L132").Value = 6630.6666
$ws.Range("M132").Value = -6392.999899999999
$ws.Range("N132").Value = -11690.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2317.875
$ws.Range("I136").Value = 1295
$ws.Range("J136").Value = 3048.5
$ws.Range("K136").Value = 3885
$ws.Range("L136").Value = 9145.5
$ws.Range("M136").Value = -1335
$ws.Range("N136").Value = -14245.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 335156.84
$ws.Range("I105").Value = 335660
$ws.Range("K105").Value = 335660
$ws.Range("M105").Value = -333913

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1948.6863
$ws.Range("I134").Value = 1837.1111
$ws.Range("K134").Value = 5511.3333
$ws.Range("M134").Value = -2976.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1624.3
$ws.Range("I31").Value = 857.0492
$ws.Range("J31").Value = 2824.359
$ws.Range("K31").Value = 857.0492
$ws.Range("L31").Value = 2824.359
$ws.Range("M31").Value = -562.0492
$ws.Range("N31").Value = -3414.359

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1624.3
$ws.Range("I34").Value = 857.0492
$ws.Range("J34").Value = 2824.359
$ws.Range("K34").Value = 857.0492
$ws.Range("L34").Value = 2824.359
$ws.Range("M34").Value = -655.0492
$ws.Range("N34").Value = -3228.359

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 6661.0435
$ws.Range("I58").Value = 1060.4546
$ws.Range("J58").Value = 20877.924
$ws.Range("K58").Value = 1060.4546
$ws.Range("L58").Value = 20877.924
$ws.Range("M58").Value = -857.4546
$ws.Range("N58").Value = -21283.924

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 3424.0571
$ws.Range("I107").Value = 6449.9414
$ws.Range("K107").Value = 6449.9414
$ws.Range("M107").Value = -4529.9414

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3812.8667
$ws.Range("I132").Value = 3620.111
$ws.Range("J132").Value = 4102
$ws.Range("K132").Value = 10860.333
$ws.Range("L132").Value = 12306
$ws.Range("M132").Value = -8330.332999999999
$ws.Range("N132").Value = -17366

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1118.5
$ws.Range("I134").Value = 1118.5
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 3355.5
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -820.5
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 6661.0435
$ws.Range("I136").Value = 1060.4546
$ws.Range("J136").Value = 20877.924
$ws.Range("K136").Value = 3181.3638
$ws.Range("L136").Value = 62633.772
$ws.Range("M136").Value = -631.3638000000001
$ws.Range("N136").Value = -67733.772

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 702808.9399999999
$ws.Range("I131").Value = 372.15384
$ws.Range("J131").Value = 814170.9
$ws.Range("K131").Value = 1116.46152
$ws.Range("L131").Value = 2442512.7
$ws.Range("M131").Value = 3923.53848
$ws.Range("N131").Value = -2452592.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 58508.137
$ws.Range("I70").Value = 90945.74000000001
$ws.Range("J70").Value = 5217.7856
$ws.Range("K70").Value = 90945.74000000001
$ws.Range("L70").Value = 5217.7856
$ws.Range("M70").Value = -90675.74000000001
$ws.Range("N70").Value = -5757.7856

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 58508.137
$ws.Range("I73").Value = 90945.74000000001
$ws.Range("J73").Value = 5217.7856
$ws.Range("K73").Value = 90945.74000000001
$ws.Range("L73").Value = 5217.7856
$ws.Range("M73").Value = -90009.74000000001
$ws.Range("N73").Value = -7089.7856

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 531920.0600000001
$ws.Range("I107").Value = 267.41666
$ws.Range("J107").Value = 1443324.6
$ws.Range("K107").Value = 267.41666
$ws.Range("L107").Value = 1443324.6
$ws.Range("M107").Value = 1652.58334
$ws.Range("N107").Value = -1447164.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1783.25
$ws.Range("I132").Value = 1833.7778
$ws.Range("J132").Value = 1631.6666
$ws.Range("K132").Value = 5501.3334
$ws.Range("L132").Value = 4894.9998
$ws.Range("M132").Value = -2971.3334
$ws.Range("N132").Value = -9954.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 858.4167
$ws.Range("I55").Value = 190
$ws.Range("J55").Value = 1283.7727
$ws.Range("K55").Value = 190
$ws.Range("L55").Value = 1283.7727
$ws.Range("M55").Value = -17
$ws.Range("N55").Value = -1629.7727

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1705.9231
$ws.Range("J93").Value = 1683.5
$ws.Range("L93").Value = 1683.5
$ws.Range("N93").Value = -4179.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H119").Value = 36661.668
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 36661.668
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 36661.668
$ws.Range("M119").ClearContents()
$ws.Range("N119").Value = -46337.668

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 43240
$ws.Range("J119").Value = 43240
$ws.Range("L119").Value = 43240
$ws.Range("N119").Value = -52916

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1404.0526
$ws.Range("I136").Value = 533.82355
$ws.Range("K136").Value = 1601.47065
$ws.Range("M136").Value = 948.5293500000002
